$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H19").Value = 9234627
$ws.Range("I19").Value = 7224515.5
$ws.Range("J19").Value = 12501059
$ws.Range("K19").Value = 7224515.5
$ws.Range("L19").Value = 12501059
$ws.Range("M19").Value = -7224340.5
$ws.Range("N19").Value = -12501409

$ws.Range("H137").Value = 6437670
$ws.Range("I137").Value = 10418172
$ws.Range("J137").Value = 68867.13
$ws.Range("K137").Value = 31254516
$ws.Range("L137").Value = 206601.39
$ws.Range("M137").Value = -31251966
$ws.Range("N137").Value = -211701.39

# ---------------------------------------------------------------------------
# Sheet CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 2037.9246
$ws.Range("I31").Value = 1617.4706
$ws.Range("J31").Value = 2236.4722
$ws.Range("K31").Value = 1617.4706
$ws.Range("L31").Value = 2236.4722
$ws.Range("M31").Value = -1322.4706
$ws.Range("N31").Value = -2826.4722

$ws.Range("H34").Value = 2037.9246
$ws.Range("I34").Value = 1617.4706
$ws.Range("J34").Value = 2236.4722
$ws.Range("K34").Value = 1617.4706
$ws.Range("L34").Value = 2236.4722
$ws.Range("M34").Value = -1415.4706
$ws.Range("N34").Value = -2640.4722

# ---------------------------------------------------------------------------
# Sheet CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H43").Value = 15000
$ws.Range("J43").Value = 15000
$ws.Range("L43").Value = 45000
$ws.Range("N43").Value = -45228

$ws.Range("H68").Value = 1692.8948
$ws.Range("I68").Value = 1255.3846
$ws.Range("J68").Value = 2059.8386
$ws.Range("K68").Value = 3766.1538
$ws.Range("L68").Value = 6179.5158
$ws.Range("M68").Value = -2955.1538
$ws.Range("N68").Value = -7801.5158

$ws.Range("H71").Value = 1692.8948
$ws.Range("I71").Value = 1255.3846
$ws.Range("J71").Value = 2059.8386
$ws.Range("K71").Value = 11298.4614
$ws.Range("L71").Value = 18538.5474
$ws.Range("M71").Value = -7242.4614
$ws.Range("N71").Value = -26650.5474

# ---------------------------------------------------------------------------
# Sheet LTW - populate previously-empty H:N columns for several leve rows
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H124").Value = 50429
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 50429
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 50429
$ws.Range("N124").Value = -60249

$ws.Range("H125").Value = 49905
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 49905
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 49905
$ws.Range("N125").Value = -59745

$ws.Range("H127").Value = 39333.332
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 39333.332
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 39333.332
$ws.Range("N127").Value = -49253.332

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0

$ws.Range("H129").Value = 20000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 20000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 20000
$ws.Range("N129").Value = -30000

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0

$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0

$ws.Range("H132").Value = 5681.421
$ws.Range("I132").Value = 4380.6924
$ws.Range("J132").Value = 8499.666999999999
$ws.Range("K132").Value = 13142.0772
$ws.Range("L132").Value = 25499.001
$ws.Range("M132").Value = -10612.0772
$ws.Range("N132").Value = -30559.001

$ws.Range("H133").Value = 55000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 55000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 55000
$ws.Range("N133").Value = -60060

$ws.Range("H134").Value = 33414.285
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 33414.285
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 33414.285
$ws.Range("N134").Value = -43554.285

$ws.Range("H135").Value = 28429
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 28429
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 28429
$ws.Range("N135").Value = -38569

$ws.Range("H136").Value = 1293.0968
$ws.Range("I136").Value = 847.6667
$ws.Range("J136").Value = 4299.75
$ws.Range("K136").Value = 2543.0001
$ws.Range("L136").Value = 12899.25
$ws.Range("M136").Value = 6.999899999999798
$ws.Range("N136").Value = -17999.25

$ws.Range("H137").Value = 33925.332
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 33925.332
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 33925.332
$ws.Range("N137").Value = -44125.332

$ws.Range("H138").Value = 16490
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 16490
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 16490
$ws.Range("N138").Value = -26770

$ws.Range("H139").Value = 38000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 38000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 38000
$ws.Range("N139").Value = -48280

$ws.Range("H140").Value = 32385.428
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 32385.428
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 32385.428
$ws.Range("N140").Value = -42745.428

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0

# ---------------------------------------------------------------------------
# Sheet WVR - clear out H:N columns for rows 119-141 (these become blank)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

for ($r = 119; $r -le 141; $r++) {
    $rng = "H" + $r + ":N" + $r
    $ws.Range($rng).ClearContents()
}
